$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 values: A=Principle, B=Start Principle, C=BuyPrice, D=SellPrice,
# E=IsShortSell, F=Price Change %, G=Date, H=Profitable
$ws.Range("A3").Value = 9945
$ws.Range("C3").Value = 110.77
$ws.Range("D3").Value = 110.16
$ws.Range("F3").Value = -0.55

# Date serial value (days since 1899-12-30), keep existing date formatting/style
$ws.Range("G3").Value = 42608.639108796298

$ws.Range("H3").Value = $false
